# Actualización 11 de Mayo - Tarde
# Updates the "Blancos" worksheet: corrects the subject (col E) / teacher
# (col F) pairing for several rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blancos")

# Map of row -> (new subject text, new teacher text)
$updates = @{
    2  = @("INGLÉS IV", "González Nuñez Veronica")
    3  = @("CÁLCULO DIFERENCIAL", "Muñoz Rivadeneyra Salvador")
    4  = @("FÍSICA I", "González Sánchez Rene Aurelio")
    5  = @("REALIZA ANÁLISIS CITOQUÍMICOS A LÍQUIDOS Y SECRECIONES CORPORALES", "Ángel Martínez Noe Cristobal")
    6  = @("REALIZA ANÁLISIS INMUNOLÓGICOS", "Rivera Serra Alma Lilian")
    7  = @("CÁLCULO DIFERENCIAL", "Muñoz Rivadeneyra Salvador")
    8  = @("INGLÉS IV", "González Nuñez Veronica")
    9  = @("FÍSICA I", "González Sánchez Rene Aurelio")
    10 = @("REALIZA ANÁLISIS HEMATOLÓGICOS DE SERIE ROJA", "Rivera Serra Alma Lilian")
    25 = @("INGLÉS IV", "González Nuñez Veronica")
    26 = @("CÁLCULO DIFERENCIAL", "Muñoz Rivadeneyra Salvador")
    28 = @("REALIZA ANÁLISIS CITOQUÍMICOS A LÍQUIDOS Y SECRECIONES CORPORALES", "Ángel Martínez Noe Cristobal")
    29 = @("REALIZA ANÁLISIS HEMATOLÓGICOS DE SERIE ROJA", "Rivera Serra Alma Lilian")
    30 = @("REALIZA ANÁLISIS INMUNOLÓGICOS", "Rivera Serra Alma Lilian")
    31 = @("INGLÉS IV", "González Nuñez Veronica")
    33 = @("FÍSICA I", "González Sánchez Rene Aurelio")
    35 = @("REALIZA ANÁLISIS INMUNOLÓGICOS", "Rivera Serra Alma Lilian")
    36 = @("CÁLCULO DIFERENCIAL", "Muñoz Rivadeneyra Salvador")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $ws.Cells.Item($row, 5).Value = $pair[0]
    $ws.Cells.Item($row, 6).Value = $pair[1]
}
